$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers to support plans without sticky-ids or levels:
# "Unique Sticky ID" -> "Row ID", "Task Name" -> "Task",
# "Start" -> "Start Date", "Finish" -> "End Date"
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update the saved selection/active cell
$ws.Range("F2").Select()
